# Auto-generated edit script - updates the "18-12-2022" coin snapshot
# table (rows 2-51) to the next GitHub Actions refresh (hour 8 -> 9,
# refreshed prices, and the bottom-of-ranking coins shifting up/down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data cells in this sheet are stored as text, even the numeric-
# looking ones (price/hour). Setting .Value on a numeric-looking string
# makes Excel coerce it to a real number, so force text with a leading
# apostrophe and then drop the resulting quote-prefix style so the cell
# is left exactly as plain, unstyled text (matching the rest of the sheet).
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '248.15'
Set-TextValue $ws.Range("G2") '9'

# Row 3
Set-TextValue $ws.Range("D3") '22.69'
Set-TextValue $ws.Range("G3") '9'

# Row 4
Set-TextValue $ws.Range("D4") '5.557'
Set-TextValue $ws.Range("G4") '9'

# Row 5
Set-TextValue $ws.Range("D5") '0.05623'
Set-TextValue $ws.Range("G5") '9'

# Row 6
Set-TextValue $ws.Range("D6") '3.394'
Set-TextValue $ws.Range("G6") '9'

# Row 7
Set-TextValue $ws.Range("D7") '6.473'
Set-TextValue $ws.Range("G7") '9'

# Row 8
Set-TextValue $ws.Range("G8") '9'

# Row 9
Set-TextValue $ws.Range("D9") '0.8029'
Set-TextValue $ws.Range("G9") '9'

# Row 10
Set-TextValue $ws.Range("D10") '0.1428'
Set-TextValue $ws.Range("G10") '9'

# Row 11
Set-TextValue $ws.Range("D11") '0.07305'
Set-TextValue $ws.Range("G11") '9'

# Row 12
Set-TextValue $ws.Range("D12") '0.03206'
Set-TextValue $ws.Range("G12") '9'

# Row 13
Set-TextValue $ws.Range("D13") '0.02994'
Set-TextValue $ws.Range("G13") '9'

# Row 14
Set-TextValue $ws.Range("D14") '0.09264'
Set-TextValue $ws.Range("G14") '9'

# Row 15
Set-TextValue $ws.Range("D15") '0.001673'
Set-TextValue $ws.Range("G15") '9'

# Row 16
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range("D16") '0.04692'
$ws.Range("E16").Value = '15CoinExTokenCET'
Set-TextValue $ws.Range("G16") '9'

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D17") '0.006275'
$ws.Range("E17").Value = '16TigerCashTCH'
Set-TextValue $ws.Range("G17") '9'

# Row 18
$ws.Range("B18").Value = 'BitKan'
$ws.Range("C18").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range("D18") '0.001051'
$ws.Range("E18").Value = '17BitKanKAN'
Set-TextValue $ws.Range("G18") '9'

# Row 19
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range("D19") '0.003835'
$ws.Range("E19").Value = '18HotbitTokenHTB'
Set-TextValue $ws.Range("G19") '9'

# Row 20
$ws.Range("B20").Value = 'NitroEx'
$ws.Range("C20").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue $ws.Range("D20") '0.0001501'
$ws.Range("E20").Value = '19NitroExNTX'
Set-TextValue $ws.Range("G20") '9'

# Row 21
$ws.Range("B21").Value = 'UpBots'
$ws.Range("C21").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue $ws.Range("D21") '0.0004004'
$ws.Range("E21").Value = '20UpBotsUBXT'
Set-TextValue $ws.Range("G21") '9'

# Row 22
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D22") '3.983'
$ws.Range("E22").Value = '21LEOLEO'
Set-TextValue $ws.Range("G22") '9'

# Row 23
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range("D23") '2.084'
$ws.Range("E23").Value = '22BTSETokenBTSE'
Set-TextValue $ws.Range("G23") '9'

# Row 24
$ws.Range("B24").Value = 'One'
$ws.Range("C24").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range("D24") '0.01168'
$ws.Range("E24").Value = '23OneONEBestin24h'
Set-TextValue $ws.Range("G24") '9'

# Row 25
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range("D25") '0.3290'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'
Set-TextValue $ws.Range("G25") '9'

# Row 26
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range("D26") '0.1292'
$ws.Range("E26").Value = '25ProBitTokenPROB'
Set-TextValue $ws.Range("G26") '9'

# Row 27
$ws.Range("B27").Value = 'MCDex'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range("D27") '2.586'
$ws.Range("E27").Value = '26MCDexMCB'
Set-TextValue $ws.Range("G27") '9'

# Row 28
Set-TextValue $ws.Range("G28") '9'

# Row 29
Set-TextValue $ws.Range("G29") '9'

# Row 30
Set-TextValue $ws.Range("G30") '9'

# Row 31
Set-TextValue $ws.Range("G31") '9'

# Row 32
Set-TextValue $ws.Range("G32") '9'

# Row 33
Set-TextValue $ws.Range("G33") '9'

# Row 34
Set-TextValue $ws.Range("G34") '9'

# Row 35
Set-TextValue $ws.Range("G35") '9'

# Row 36
Set-TextValue $ws.Range("G36") '9'

# Row 37
Set-TextValue $ws.Range("G37") '9'

# Row 38
Set-TextValue $ws.Range("G38") '9'

# Row 39
Set-TextValue $ws.Range("G39") '9'

# Row 40
Set-TextValue $ws.Range("D40") '0.04208'
Set-TextValue $ws.Range("G40") '9'

# Row 41
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws.Range("D41") '0.007007'
$ws.Range("E41").Value = '40KickTokenKICK'
Set-TextValue $ws.Range("G41") '9'

# Row 42
Set-TextValue $ws.Range("D42") '0.003503'
Set-TextValue $ws.Range("G42") '9'

# Row 43
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws.Range("D43") '0.1046'
$ws.Range("E43").Value = '42BKEXTokenBKK'
Set-TextValue $ws.Range("G43") '9'

# Row 44
Set-TextValue $ws.Range("D44") '0.008751'
Set-TextValue $ws.Range("G44") '9'

# Row 45
Set-TextValue $ws.Range("D45") '0.00005638'
Set-TextValue $ws.Range("G45") '9'

# Row 46
Set-TextValue $ws.Range("G46") '9'

# Row 47
Set-TextValue $ws.Range("D47") '0.6807'
Set-TextValue $ws.Range("G47") '9'

# Row 48
Set-TextValue $ws.Range("D48") '0.02698'
$ws.Range("E48").Value = '47BOLOBOLOWorstin24h'
Set-TextValue $ws.Range("G48") '9'

# Row 49
Set-TextValue $ws.Range("G49") '9'

# Row 50
Set-TextValue $ws.Range("G50") '9'

# Row 51
Set-TextValue $ws.Range("G51") '9'
